{"js": "// Apply the \"Added many more features\" edit to the Hot Hot Halloween review.\n// Each entry is an exact-text search/replace; results are matched by full\n// containing text so formatting (heading style, bold/italic runs) is kept.\nconst replacements = [\n  {\n    find: \"Play Hot Hot Halloween Slot for Free - Review and Rating\",\n    replace: \"Play Hot Hot Halloween Slot for Free\",\n  },\n  {\n    find: \"Unique and spooky horizontal game mechanic with great winning potential\",\n    replace: \"Unique horizontal game mechanic\",\n  },\n  {\n    find: \"High RTP of 96.62% and high volatility for huge payouts\",\n    replace: \"Free spins, double symbols, and wild wins\",\n  },\n  {\n    find: \"Double and triple symbols, wilds, and free spins feature for more chances to win\",\n    replace: \"High RTP of 96.62%\",\n  },\n  {\n    find: \"Playable on both PC and mobile, HTML5 technology compatible on Android or iPhone\",\n    replace: \"Compatible on both PC and mobile devices\",\n  },\n  {\n    find: \"High volatility may not be suitable for players who prefer low-risk slots\",\n    replace: \"High volatility\",\n  },\n  {\n    find: \"Bonus features are triggered randomly and can't be activated on demand\",\n    replace: \"Max bet of \\u20ac5,000.00 may be too high for some players\",\n  },\n  {\n    find: \"Ready to play Hot Hot Halloween for free? Check out our review to learn more about its unique horizontal game mechanic, great winning potential, and spooky features.\",\n    replace: \"Read our review of Hot Hot Halloween slot game and play for free. Enjoy the spooky Halloween theme and high RTP!\",\n  },\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"Added many more features\" edit to the Hot Hot Halloween review.\n# Each pair is an exact-text find/replace across the whole document body;\n# wdReplaceAll keeps run formatting (heading style, bold/italic) intact.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"Play Hot Hot Halloween Slot for Free - Review and Rating\"; Replace = \"Play Hot Hot Halloween Slot for Free\" },\n    @{ Find = \"Unique and spooky horizontal game mechanic with great winning potential\"; Replace = \"Unique horizontal game mechanic\" },\n    @{ Find = \"High RTP of 96.62% and high volatility for huge payouts\"; Replace = \"Free spins, double symbols, and wild wins\" },\n    @{ Find = \"Double and triple symbols, wilds, and free spins feature for more chances to win\"; Replace = \"High RTP of 96.62%\" },\n    @{ Find = \"Playable on both PC and mobile, HTML5 technology compatible on Android or iPhone\"; Replace = \"Compatible on both PC and mobile devices\" },\n    @{ Find = \"High volatility may not be suitable for players who prefer low-risk slots\"; Replace = \"High volatility\" },\n    @{ Find = \"Bonus features are triggered randomly and can't be activated on demand\"; Replace = \"Max bet of \u20ac5,000.00 may be too high for some players\" },\n    @{ Find = \"Ready to play Hot Hot Halloween for free? Check out our review to learn more about its unique horizontal game mechanic, great winning potential, and spooky features.\"; Replace = \"Read our review of Hot Hot Halloween slot game and play for free. Enjoy the spooky Halloween theme and high RTP!\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Find\n    $find.Replacement.Text = $r.Replace\n    $find.Execute($r.Find, $false, $false, $false, $false, $false, $true, 1, $false, $r.Replace, 2)\n}\n"}
